# Edit: add new sheet "main sistem lagi" (sheetId 11) after "main sistem baru",
# containing an expanded version of the web-scraping summarizer script (rows 1-40),
# and leave "main sistem baru" as a plain (non-active) sheet with its original selection.

$wb = $excel.ActiveWorkbook

# 1) Create the new sheet by copying "main sistem baru" so it inherits the same
#    column styling/widths (style "2" = Courier New 9pt, left/top aligned).
$src = $wb.Worksheets.Item("main sistem baru")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $lastSheet)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "main sistem lagi"

$rowValues = @{
    1 = "from IPython.display import Markdown, display"
    2 = "document2 = input(`"URL berita BBC Indonesia: `") "
    4 = "print(`"Jenis compression rate:\n 1. 5% \n 2. 10% \n 3. 20% \n 4. 30%`")"
    5 = "cr2 = input(`"Pilih jenis compression rate: `")"
    8 = "if cr2 in [`"1`",`"2`",`"3`",`"4`"] :"
    9 = "    stopword = open(`"../stopword_list_tala.txt`", `"r`")"
    10 = "    stopwords = stopword.read().split(`"\n`")"
    11 = "    document3, title = get_document(document2)"
    13 = "    if document3 == `"error`" or title == `"none`":"
    14 = "        display(Markdown('**sorry, i cant access the url**'))"
    15 = "    else:      "
    16 = "        cleaning_result2 = get_clean_corpus(corpus=document3, stopwords=stopwords)"
    17 = "        terms_frequency2, df_idf2  = get_term_weighting_score(cleaning_result=cleaning_result2)"
    18 = "        "
    19 = "        if cr2 == `"1`":"
    20 = "            crate = 0.05"
    21 = "        elif cr2 == `"2`":"
    22 = "            crate = 0.1"
    23 = "        elif cr2 == `"3`":"
    24 = "            crate = 0.2"
    25 = "        elif cr2 == `"4`":"
    26 = "            crate = 0.3"
    28 = "        percobaan2 = Graph(result_doc=cleaning_result2, raw_frequency=terms_frequency2, "
    29 = "                              idf=df_idf2, cr=float(crate))"
    30 = "        display(Markdown('**Judul berita:**'))"
    31 = "        print(title)"
    33 = "        summarize = `"`""
    34 = "        for ringkasan in percobaan2.summarize:"
    35 = "            summarize += ringkasan + `" `""
    37 = "        display(Markdown('**Hasil ringkasan:**'))"
    38 = "        print(summarize)"
    39 = "else:"
    40 = "    display(Markdown('**sorry, wrong input**'))"
}

# 2) Clear out the copied content first (column A numbering + column B text),
#    then rewrite both columns 1..40 from scratch.
$new.Range("A1:B24").ClearContents()

foreach ($r in 1..40) {
    $new.Cells.Item($r, 1).Value = $r
    if ($rowValues.ContainsKey($r)) {
        $new.Cells.Item($r, 2).Value = $rowValues[$r]
    }
}

# 3) Match the new sheet's page setup / view as closely as possible.
$new.PageSetup.Orientation = 1

$new.Range("E9").Select()

Write-Output "done"
